# Generate Report for handoff
#
# The b8f8fc47-5348-4846-b38b-cfaaf95350a0 file is now ready to be handed
# off again: its Status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the Latest Handoff Datetime for its language
# rows gets refreshed to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b8f8fc47-... (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for b8f8fc47-... (row 3), plus shared handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "2016-02-17 03:34:01"
$wsZhCn.Range("D3").Value = "2016-02-17 03:34:01"

# --- de-de sheet: row for b8f8fc47-... (row 3), plus shared handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "2016-02-17 03:34:11"
$wsDeDe.Range("D3").Value = "2016-02-17 03:34:11"
